$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Vostro / Dell" DKS enrollment-device rows (original rows 6-7, ids 589 & 638)
# are removed; everything below shifts up by two rows.
$ws.Rows("6:7").Delete()

# Restore the cursor/selection left behind in the authored file.
$ws.Range("E16").Select()

# Page setup metadata captured when the workbook was saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
